$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to reflect the new extraction run (date moved from
# 2024-11-25 09:35:07 to 2024-11-26 09:08:21)
$ws.Name = "IClientBalance-20241126-090821-"

# Column G holds "Dt. Referencia" as a date serial; every data row
# (2..274) moves one day forward: 45621 (2024-11-25) -> 45622 (2024-11-26)
for ($r = 2; $r -le 274; $r++) {
    $ws.Cells.Item($r, 7).Value = 45622
}

# Row 120: Saldo Previsto / Vl. Total corrected from 22823.18 to 12823.18
$ws.Cells.Item(120, 5).Value = 12823.18
$ws.Cells.Item(120, 8).Value = 12823.18

# Row 129: Saldo Previsto / Vl. Total corrected from 2234.33 to 0
$ws.Cells.Item(129, 5).Value = 0
$ws.Cells.Item(129, 8).Value = 0
